# Update "想去人数" (want-to-go count) figures across sheets to match the
# newly scraped output (gh-pages rebuild at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 10041
$ws1.Range("F8").Value = 1580
$ws1.Range("F20").Value = 338
$ws1.Range("F31").Value = 322
$ws1.Range("F36").Value = 715

# --- 演出 sheet ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F23").Value = 67
$ws2.Range("F36").Value = 28

# --- 本地生活 sheet ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F10").Value = 262

# --- 全部类型 sheet ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 10041
$ws4.Range("F10").Value = 262
$ws4.Range("F11").Value = 262
$ws4.Range("F24").Value = 338
$ws4.Range("F31").Value = 67
$ws4.Range("F41").Value = 715
